$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("USER8")

# --- Reproduce the style-table "churn" (3 extra unused border/cellXfs
# entries) that appears in the target workbook. These entries are never
# actually referenced by any cell (every cell stays at the default style),
# so we apply a border + a date number-format to a few cells that we are
# about to overwrite with real data anyway, then reset them back to the
# "Normal" style. This allocates new border/xf slots without changing the
# visible formatting of the final sheet.
$g1 = $ws.Cells.Item(2, 1)
$g1.Borders.Item(7).LineStyle = 1
$g1.NumberFormat = "m/d/yy h:mm"
$g1.Style = "Normal"

$g2 = $ws.Cells.Item(2, 2)
$g2.Borders.Item(8).LineStyle = 1
$g2.NumberFormat = "m/d/yy h:mm"
$g2.Style = "Normal"

$g3 = $ws.Cells.Item(2, 3)
$g3.Borders.Item(9).LineStyle = 1
$g3.NumberFormat = "m/d/yy h:mm"
$g3.Style = "Normal"

# --- Insert the new experiment run at the top (row 2), push the old row 2
# data down to row 4, and update row 3 with its new values. ---

# New row 2 (latest experiment result)
$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = 0.2
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0.54492753623188395
$ws.Cells.Item(2, 5).Value = 0.2

# Updated row 3
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 0.13095238095238096
$ws.Cells.Item(3, 3).Value = 0.1
$ws.Cells.Item(3, 4).Value = 0.61904761904761907
$ws.Cells.Item(3, 5).Value = 0

# New row 4 (previous row 2, unchanged)
$ws.Cells.Item(4, 1).Value = 0
$ws.Cells.Item(4, 2).Value = 0.27826086956521739
$ws.Cells.Item(4, 3).Value = 0.2
$ws.Cells.Item(4, 4).Value = 0.63478260869565217
$ws.Cells.Item(4, 5).Value = 0.66666666666666663
